$wb = $excel.ActiveWorkbook

# --- first_query ---
$ws = $wb.Worksheets.Item("first_query")
$arr = New-Object 'object[,]' 31,4
$arr[0,0] = 107
$arr[0,1] = 95
$arr[0,2] = 113
$arr[0,3] = 91
$arr[1,0] = 65
$arr[1,1] = 61
$arr[1,2] = 65
$arr[1,3] = 59
$arr[2,0] = 65
$arr[2,1] = 62
$arr[2,2] = 65
$arr[2,3] = 62
$arr[3,0] = 64
$arr[3,1] = 60
$arr[3,2] = 65
$arr[3,3] = 59
$arr[4,0] = 63
$arr[4,1] = 62
$arr[4,2] = 71
$arr[4,3] = 62
$arr[5,0] = 64
$arr[5,1] = 63
$arr[5,2] = 67
$arr[5,3] = 60
$arr[6,0] = 62
$arr[6,1] = 63
$arr[6,2] = 64
$arr[6,3] = 62
$arr[7,0] = 67
$arr[7,1] = 61
$arr[7,2] = 63
$arr[7,3] = 62
$arr[8,0] = 62
$arr[8,1] = 62
$arr[8,2] = 64
$arr[8,3] = 65
$arr[9,0] = 62
$arr[9,1] = 62
$arr[9,2] = 64
$arr[9,3] = 62
$arr[10,0] = 61
$arr[10,1] = 61
$arr[10,2] = 64
$arr[10,3] = 61
$arr[11,0] = 63
$arr[11,1] = 66
$arr[11,2] = 69
$arr[11,3] = 60
$arr[12,0] = 63
$arr[12,1] = 61
$arr[12,2] = 64
$arr[12,3] = 59
$arr[13,0] = 62
$arr[13,1] = 64
$arr[13,2] = 63
$arr[13,3] = 62
$arr[14,0] = 62
$arr[14,1] = 64
$arr[14,2] = 63
$arr[14,3] = 60
$arr[15,0] = 62
$arr[15,1] = 64
$arr[15,2] = 61
$arr[15,3] = 61
$arr[16,0] = 62
$arr[16,1] = 61
$arr[16,2] = 60
$arr[16,3] = 62
$arr[17,0] = 65
$arr[17,1] = 64
$arr[17,2] = 63
$arr[17,3] = 62
$arr[18,0] = 68
$arr[18,1] = 74
$arr[18,2] = 62
$arr[18,3] = 62
$arr[19,0] = 69
$arr[19,1] = 78
$arr[19,2] = 63
$arr[19,3] = 62
$arr[20,0] = 68
$arr[20,1] = 65
$arr[20,2] = 62
$arr[20,3] = 60
$arr[21,0] = 63
$arr[21,1] = 68
$arr[21,2] = 63
$arr[21,3] = 62
$arr[22,0] = 62
$arr[22,1] = 67
$arr[22,2] = 62
$arr[22,3] = 62
$arr[23,0] = 75
$arr[23,1] = 64
$arr[23,2] = 62
$arr[23,3] = 61
$arr[24,0] = 65
$arr[24,1] = 63
$arr[24,2] = 209
$arr[24,3] = 61
$arr[25,0] = 66
$arr[25,1] = 62
$arr[25,2] = 64
$arr[25,3] = 61
$arr[26,0] = 67
$arr[26,1] = 65
$arr[26,2] = 59
$arr[26,3] = 64
$arr[27,0] = 64
$arr[27,1] = 62
$arr[27,2] = 62
$arr[27,3] = 61
$arr[28,0] = 68
$arr[28,1] = 62
$arr[28,2] = 62
$arr[28,3] = 62
$arr[29,0] = 64
$arr[29,1] = 64
$arr[29,2] = 62
$arr[29,3] = 62
$arr[30,0] = 65
$arr[30,1] = 64
$arr[30,2] = 61
$arr[30,3] = 62
$ws.Range("A2:D32").Value = $arr

# --- second_query ---
$ws = $wb.Worksheets.Item("second_query")
$arr = New-Object 'object[,]' 31,4
$arr[0,0] = 78
$arr[0,1] = 79
$arr[0,2] = 76
$arr[0,3] = 83
$arr[1,0] = 53
$arr[1,1] = 53
$arr[1,2] = 57
$arr[1,3] = 49
$arr[2,0] = 54
$arr[2,1] = 55
$arr[2,2] = 50
$arr[2,3] = 51
$arr[3,0] = 51
$arr[3,1] = 54
$arr[3,2] = 51
$arr[3,3] = 52
$arr[4,0] = 50
$arr[4,1] = 49
$arr[4,2] = 53
$arr[4,3] = 51
$arr[5,0] = 50
$arr[5,1] = 52
$arr[5,2] = 52
$arr[5,3] = 51
$arr[6,0] = 51
$arr[6,1] = 51
$arr[6,2] = 51
$arr[6,3] = 51
$arr[7,0] = 49
$arr[7,1] = 52
$arr[7,2] = 53
$arr[7,3] = 52
$arr[8,0] = 56
$arr[8,1] = 50
$arr[8,2] = 50
$arr[8,3] = 50
$arr[9,0] = 51
$arr[9,1] = 50
$arr[9,2] = 52
$arr[9,3] = 51
$arr[10,0] = 51
$arr[10,1] = 50
$arr[10,2] = 52
$arr[10,3] = 50
$arr[11,0] = 52
$arr[11,1] = 50
$arr[11,2] = 54
$arr[11,3] = 52
$arr[12,0] = 50
$arr[12,1] = 55
$arr[12,2] = 50
$arr[12,3] = 51
$arr[13,0] = 50
$arr[13,1] = 57
$arr[13,2] = 50
$arr[13,3] = 173
$arr[14,0] = 51
$arr[14,1] = 51
$arr[14,2] = 50
$arr[14,3] = 48
$arr[15,0] = 50
$arr[15,1] = 53
$arr[15,2] = 48
$arr[15,3] = 52
$arr[16,0] = 51
$arr[16,1] = 50
$arr[16,2] = 52
$arr[16,3] = 50
$arr[17,0] = 48
$arr[17,1] = 50
$arr[17,2] = 49
$arr[17,3] = 50
$arr[18,0] = 60
$arr[18,1] = 80
$arr[18,2] = 50
$arr[18,3] = 51
$arr[19,0] = 55
$arr[19,1] = 51
$arr[19,2] = 49
$arr[19,3] = 53
$arr[20,0] = 51
$arr[20,1] = 53
$arr[20,2] = 50
$arr[20,3] = 54
$arr[21,0] = 50
$arr[21,1] = 52
$arr[21,2] = 48
$arr[21,3] = 50
$arr[22,0] = 51
$arr[22,1] = 53
$arr[22,2] = 49
$arr[22,3] = 50
$arr[23,0] = 50
$arr[23,1] = 51
$arr[23,2] = 50
$arr[23,3] = 50
$arr[24,0] = 262
$arr[24,1] = 51
$arr[24,2] = 54
$arr[24,3] = 51
$arr[25,0] = 60
$arr[25,1] = 56
$arr[25,2] = 49
$arr[25,3] = 50
$arr[26,0] = 55
$arr[26,1] = 52
$arr[26,2] = 47
$arr[26,3] = 52
$arr[27,0] = 54
$arr[27,1] = 51
$arr[27,2] = 52
$arr[27,3] = 50
$arr[28,0] = 55
$arr[28,1] = 52
$arr[28,2] = 49
$arr[28,3] = 53
$arr[29,0] = 53
$arr[29,1] = 51
$arr[29,2] = 50
$arr[29,3] = 50
$arr[30,0] = 54
$arr[30,1] = 52
$arr[30,2] = 49
$arr[30,3] = 49
$ws.Range("A2:D32").Value = $arr

# --- third_query ---
$ws = $wb.Worksheets.Item("third_query")
$arr = New-Object 'object[,]' 31,4
$arr[0,0] = 64
$arr[0,1] = 64
$arr[0,2] = 62
$arr[0,3] = 59
$arr[1,0] = 41
$arr[1,1] = 38
$arr[1,2] = 37
$arr[1,3] = 41
$arr[2,0] = 38
$arr[2,1] = 37
$arr[2,2] = 41
$arr[2,3] = 37
$arr[3,0] = 37
$arr[3,1] = 35
$arr[3,2] = 41
$arr[3,3] = 36
$arr[4,0] = 37
$arr[4,1] = 41
$arr[4,2] = 42
$arr[4,3] = 36
$arr[5,0] = 38
$arr[5,1] = 36
$arr[5,2] = 40
$arr[5,3] = 35
$arr[6,0] = 37
$arr[6,1] = 36
$arr[6,2] = 42
$arr[6,3] = 37
$arr[7,0] = 37
$arr[7,1] = 37
$arr[7,2] = 37
$arr[7,3] = 36
$arr[8,0] = 37
$arr[8,1] = 36
$arr[8,2] = 38
$arr[8,3] = 36
$arr[9,0] = 37
$arr[9,1] = 36
$arr[9,2] = 39
$arr[9,3] = 37
$arr[10,0] = 36
$arr[10,1] = 36
$arr[10,2] = 38
$arr[10,3] = 37
$arr[11,0] = 36
$arr[11,1] = 177
$arr[11,2] = 39
$arr[11,3] = 39
$arr[12,0] = 37
$arr[12,1] = 36
$arr[12,2] = 36
$arr[12,3] = 38
$arr[13,0] = 38
$arr[13,1] = 39
$arr[13,2] = 38
$arr[13,3] = 41
$arr[14,0] = 38
$arr[14,1] = 39
$arr[14,2] = 36
$arr[14,3] = 36
$arr[15,0] = 37
$arr[15,1] = 38
$arr[15,2] = 36
$arr[15,3] = 38
$arr[16,0] = 36
$arr[16,1] = 37
$arr[16,2] = 36
$arr[16,3] = 36
$arr[17,0] = 37
$arr[17,1] = 40
$arr[17,2] = 36
$arr[17,3] = 36
$arr[18,0] = 38
$arr[18,1] = 54
$arr[18,2] = 36
$arr[18,3] = 36
$arr[19,0] = 37
$arr[19,1] = 36
$arr[19,2] = 36
$arr[19,3] = 37
$arr[20,0] = 35
$arr[20,1] = 37
$arr[20,2] = 36
$arr[20,3] = 34
$arr[21,0] = 36
$arr[21,1] = 38
$arr[21,2] = 39
$arr[21,3] = 37
$arr[22,0] = 38
$arr[22,1] = 40
$arr[22,2] = 37
$arr[22,3] = 36
$arr[23,0] = 38
$arr[23,1] = 37
$arr[23,2] = 36
$arr[23,3] = 37
$arr[24,0] = 36
$arr[24,1] = 36
$arr[24,2] = 36
$arr[24,3] = 36
$arr[25,0] = 36
$arr[25,1] = 39
$arr[25,2] = 36
$arr[25,3] = 152
$arr[26,0] = 37
$arr[26,1] = 57
$arr[26,2] = 35
$arr[26,3] = 36
$arr[27,0] = 36
$arr[27,1] = 37
$arr[27,2] = 37
$arr[27,3] = 37
$arr[28,0] = 43
$arr[28,1] = 37
$arr[28,2] = 37
$arr[28,3] = 37
$arr[29,0] = 39
$arr[29,1] = 38
$arr[29,2] = 36
$arr[29,3] = 36
$arr[30,0] = 38
$arr[30,1] = 41
$arr[30,2] = 36
$arr[30,3] = 38
$ws.Range("A2:D32").Value = $arr

# --- fourth_query ---
$ws = $wb.Worksheets.Item("fourth_query")
$arr = New-Object 'object[,]' 31,4
$arr[0,0] = 77
$arr[0,1] = 117
$arr[0,2] = 152
$arr[0,3] = 163
$arr[1,0] = 26
$arr[1,1] = 47
$arr[1,2] = 73
$arr[1,3] = 84
$arr[2,0] = 24
$arr[2,1] = 47
$arr[2,2] = 76
$arr[2,3] = 80
$arr[3,0] = 25
$arr[3,1] = 49
$arr[3,2] = 68
$arr[3,3] = 80
$arr[4,0] = 25
$arr[4,1] = 46
$arr[4,2] = 72
$arr[4,3] = 79
$arr[5,0] = 25
$arr[5,1] = 46
$arr[5,2] = 75
$arr[5,3] = 84
$arr[6,0] = 24
$arr[6,1] = 45
$arr[6,2] = 67
$arr[6,3] = 82
$arr[7,0] = 24
$arr[7,1] = 45
$arr[7,2] = 203
$arr[7,3] = 79
$arr[8,0] = 25
$arr[8,1] = 45
$arr[8,2] = 69
$arr[8,3] = 81
$arr[9,0] = 24
$arr[9,1] = 47
$arr[9,2] = 78
$arr[9,3] = 80
$arr[10,0] = 30
$arr[10,1] = 48
$arr[10,2] = 70
$arr[10,3] = 79
$arr[11,0] = 29
$arr[11,1] = 46
$arr[11,2] = 70
$arr[11,3] = 78
$arr[12,0] = 24
$arr[12,1] = 50
$arr[12,2] = 67
$arr[12,3] = 81
$arr[13,0] = 24
$arr[13,1] = 48
$arr[13,2] = 69
$arr[13,3] = 79
$arr[14,0] = 24
$arr[14,1] = 47
$arr[14,2] = 68
$arr[14,3] = 81
$arr[15,0] = 25
$arr[15,1] = 46
$arr[15,2] = 69
$arr[15,3] = 83
$arr[16,0] = 24
$arr[16,1] = 45
$arr[16,2] = 68
$arr[16,3] = 81
$arr[17,0] = 25
$arr[17,1] = 54
$arr[17,2] = 66
$arr[17,3] = 80
$arr[18,0] = 25
$arr[18,1] = 57
$arr[18,2] = 66
$arr[18,3] = 80
$arr[19,0] = 25
$arr[19,1] = 56
$arr[19,2] = 66
$arr[19,3] = 81
$arr[20,0] = 26
$arr[20,1] = 50
$arr[20,2] = 67
$arr[20,3] = 83
$arr[21,0] = 23
$arr[21,1] = 49
$arr[21,2] = 66
$arr[21,3] = 82
$arr[22,0] = 25
$arr[22,1] = 51
$arr[22,2] = 67
$arr[22,3] = 81
$arr[23,0] = 24
$arr[23,1] = 47
$arr[23,2] = 67
$arr[23,3] = 81
$arr[24,0] = 28
$arr[24,1] = 48
$arr[24,2] = 70
$arr[24,3] = 81
$arr[25,0] = 26
$arr[25,1] = 50
$arr[25,2] = 64
$arr[25,3] = 79
$arr[26,0] = 25
$arr[26,1] = 46
$arr[26,2] = 67
$arr[26,3] = 80
$arr[27,0] = 44
$arr[27,1] = 46
$arr[27,2] = 67
$arr[27,3] = 80
$arr[28,0] = 28
$arr[28,1] = 46
$arr[28,2] = 66
$arr[28,3] = 79
$arr[29,0] = 25
$arr[29,1] = 46
$arr[29,2] = 67
$arr[29,3] = 82
$arr[30,0] = 44
$arr[30,1] = 47
$arr[30,2] = 66
$arr[30,3] = 82
$ws.Range("A2:D32").Value = $arr

Write-Host "Update complete"